$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume data cells retain text formatting (values contain
# thousands separators as literal dots, e.g. "29.155.06") so Excel
# does not reinterpret them as numbers/dates. Restricted to the data
# rows (2-51) so the header row's existing style is left untouched.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.155.06"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "1.852.81"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "237.37"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "0.6851"
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "0.07635"
$ws.Range("E8").Value = "  +3.91%  "
$ws.Range("D9").Value = "0.3039"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "23.12"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.870.32"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7223"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "5.180"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "89.38"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "29.161.23"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "0.000007802"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "5.714"
$ws.Range("E18").Value = "  -5.03%  "
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "233.58"
$ws.Range("E20").Value = "  -5.35%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "2.099.42"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "7.415"
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("D25").Value = "161.33"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "8.940"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "0.1425"
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").Value = "18.01"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").Value = "1.953"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "1.396"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "4.505"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "4.005"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").Value = "0.05146"
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("D35").Value = "1.181"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").Value = "0.7034"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "1.021"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "2.671"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").Value = "0.01847"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "2.678"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").Value = "0.9055"
$ws.Range("D42").Value = "1.104.44"
$ws.Range("E42").Value = "  +6.25%  "
$ws.Range("D43").Value = "5.948"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").Value = "0.4275"
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("D45").Value = "69.92"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "102.18"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").Value = "1.773"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").Value = "1.996.10"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").Value = "9.131"
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("D51").Value = "6.939"
$ws.Range("E51").Value = "  -6.99%  "
